# Append new Customer/Account rows (9-16) to the "LCY Current Accounts" sheet.
#
# The new TC/Account_ID/Customer_ID/PD values all happen to look numeric
# (e.g. "118451", "1008617465" ...), so a plain Range.Value assignment would
# have Excel auto-convert them to numbers, and forcing text via a leading
# apostrophe / NumberFormat="@" would stamp a new (quotePrefix) cell style on
# them - neither matches the source data, which stores them as plain shared
# strings with the default style (just like the existing rows).
#
# Trick: write each value as a formula that evaluates to that literal text
# (="118451"), then Copy + PasteSpecial(values-only) over itself. That
# replaces the formula with its static text result while keeping the cell's
# original (default) style untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @{ Row = 9;  Values = @("118451", "1008617465", "17704777", "1001") },
    @{ Row = 10; Values = @("118451", "1008617516", "17704794", "1001") },
    @{ Row = 11; Values = @("118451", "1008617534", "17704802", "1001") },
    @{ Row = 12; Values = @("118451", "1008617535", "17704803", "1001") },
    @{ Row = 13; Values = @("118451", "1008617545", "17704810", "1001") },
    @{ Row = 14; Values = @("118451", "1008617546", "17704811", "1001") },
    @{ Row = 15; Values = @("118452", "1008617549", "17704816", "1150") },
    @{ Row = 16; Values = @("118452", "1008617551", "17704818", "1005") }
)

$xlPasteValues = -4163

foreach ($entry in $rows) {
    $r = $entry.Row
    $vals = $entry.Values

    $rowRange = $ws.Range("A" + $r + ":D" + $r)

    $ws.Range("A" + $r).Formula = '="' + $vals[0] + '"'
    $ws.Range("B" + $r).Formula = '="' + $vals[1] + '"'
    $ws.Range("C" + $r).Formula = '="' + $vals[2] + '"'
    $ws.Range("D" + $r).Formula = '="' + $vals[3] + '"'

    $rowRange.Copy()
    $rowRange.PasteSpecial($xlPasteValues)
}

$excel.CutCopyMode = $false
